# Added few more teams entry in input file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "fifthgithubrepo" team (rows 12-16) to "production".
# Write column-by-column (A for all rows, then B for all rows) so shared
# strings are interned in the same order the source workbook used.
for ($r = 12; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "gk-aks-Digital/production"
}
for ($r = 12; $r -le 16; $r++) {
    $ws.Cells.Item($r, 2).Value = "production"
}

# Append five new rows (17-21) for the additional teams under "production".
# columns: A, B, C, D, E, F, G, H
$newRows = @(
    @("gk-aks-Digital/production", "production", "my-account",          "gk-aks-Digital/my-account",          "CONFIDENTIAL", "gk-aks-CONFIDENTIAL", "my-account-write",         "gk-aks-CONFIDENTIAL\my-account-write"),
    @("gk-aks-Digital/production", "production", "cloud-foundation",     "gk-aks-Digital/cloud-foundation",    "CONFIDENTIAL", "gk-aks-CONFIDENTIAL", "cloud-foundation-write",   "gk-aks-CONFIDENTIAL\cloud-foundation-write"),
    @("gk-aks-Digital/production", "production", "dotcom",               "gk-aks-Digital/dotcom",              "CONFIDENTIAL", "gk-aks-CONFIDENTIAL", "dotcom-write",             "gk-aks-CONFIDENTIAL\dotcom-write"),
    @("gk-aks-Digital/production", "production", "platform",             "gk-aks-Digital/platform",            "CONFIDENTIAL", "gk-aks-CONFIDENTIAL", "platform-write",           "gk-aks-CONFIDENTIAL\platform-write"),
    @("gk-aks-Digital/production", "production", "front-end-platform",   "gk-aks-Digital/front-end-platform",  "CONFIDENTIAL", "gk-aks-CONFIDENTIAL", "frontend-platform-write",  "gk-aks-CONFIDENTIAL\frontend-platform-write")
)

$firstRow = 17
for ($c = 1; $c -le 8; $c++) {
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $ws.Cells.Item($firstRow + $i, $c).Value = $newRows[$i][$c - 1]
    }
}

# Update the view: scroll so column D is left-most visible, select H19
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("H19").Select()
